$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text number format on D and E columns for changed cells so that
# numeric-looking strings (e.g. "0.172", "1.00", "67.137.19") are preserved
# as literal text rather than being auto-converted to numbers by Excel.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.137.19'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +4.52%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.468.97'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +4.41%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '585.34'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +6.03%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '186.67'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +7.91%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.464.71'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +4.54%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.172'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.649'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +2.64%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '56.22'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +6.16%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.43'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +4.21%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.029.63'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +4.47%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '18.80'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +4.08%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.473.80'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +4.47%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '67.151.42'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +4.43%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +4.37%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +3.63%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '487.29'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +7.85%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.33'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +6.48%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '16.91'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +21.78%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.48'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +10.70%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '89.90'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +2.98%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.95'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +3.35%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.96'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +4.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.14'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +6.79%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '31.37'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.75%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.21'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +11.46%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '601.23'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +5.32%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.73'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +3.25%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '63.81'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +2.47%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +5.18%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.149'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +4.94%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '36.58'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +3.91%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.54'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.385'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +5.67%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0₃0759'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +4.74%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.257.76'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +6.32%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.90'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +6.47%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +4.25%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.53'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +4.13%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.79'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +24.84%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.27'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +11.95%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.77'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +7.12%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.07%  '
